$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the hidden defined name used by the "MySQL for Excel" add-in
#    (LOCAL_MYSQL_DATE_FORMAT) as seen in the new workbook.xml.
$dateFormatName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$dateFormatName.Visible = $false

# 2. Add the new "Icon" field row (row 19) under the existing table.
#    Set the shared-string cells in the same order they first appear in
#    the target file so the new shared-string entries land on the
#    expected indices (47 = "显示图标", 48 = "Icon").
$ws.Range("J19").Value = "显示图标"
$ws.Range("A19").Value = "Icon"
$ws.Range("B19").Value = "string"
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = $false
$ws.Range("F19").Value = $true
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = "Friend"

# Match the text number format ("@") used by the other cells in columns
# A, B, I and J so the new row gets the same style index as row 18.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("J19").NumberFormat = "@"

# 3. Update the active selection to match the post-edit workbook state.
[void]$ws.Range("H24").Select()

# 4. Re-create the boolean list data validation so it spans the whole
#    column below the header (F2:F1048576) instead of the split range
#    "F20:F1048576 F2:F18".
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
